# chore: adapt column header formatting to respective input file names
# Renames the "_old"/"_new" column header suffixes to the respective
# formatversion names ("_FV2410" / "_FV2504"), turns the data range into a
# proper Excel Table (so the new headers are backed by a ListObject), and
# freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row -------------------------------------------------
# Columns A-J were "<name>_old", K is "diff", L-U were "<name>_new".
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2410"
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2504"
}

# --- 2. Turn A1:U62 into an Excel Table -----------------------------------
$tableRange = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
